$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.631.07'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '3.511.39'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '609.90'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '152.29'
$ws.Range("E6").Value = '  +1.40%  '
$ws.Range("D7").Value = '3.511.09'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("D11").Value = '7.62'
$ws.Range("E11").Value = '  +8.26%  '
$ws.Range("E12").Value = '  +1.80%  '
$ws.Range("E13").Value = '  +2.59%  '
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = '4.104.23'
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").Value = '3.518.37'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '67.520.15'
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("E19").Value = '  +2.60%  '
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("E21").Value = '  +6.76%  '
$ws.Range("D22").Value = '447.06'
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("E23").Value = '  +1.23%  '
$ws.Range("E24").Value = '  +1.29%  '
$ws.Range("D25").Value = '3.651.40'
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -0.99%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '10.12'
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '8.79'
$ws.Range("E29").Value = '  +5.36%  '
$ws.Range("D30").Value = '2.52'
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("D32").Value = '0.171'
$ws.Range("E32").Value = '  +4.44%  '
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("E36").Value = '  +1.91%  '
$ws.Range("D37").Value = '3.506.80'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '8.03'
$ws.Range("E38").Value = '  +0.50%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  +7.94%  '
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("E42").Value = '  +2.91%  '
$ws.Range("D43").Value = '173.30'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").Value = '5.46'
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("D45").Value = '30.26'
$ws.Range("E45").Value = '  +9.28%  '
$ws.Range("D46").Value = '0.886'
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("D47").Value = '46.68'
$ws.Range("E47").Value = '  +2.48%  '
$ws.Range("E48").Value = '  +4.04%  '
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("D50").Value = '7.65'
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("E51").Value = '  +0.04%  '
